# Rename the sheet from "Sheet2" to "Employees"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Employees"

# Write the employee record header/data row (A1:F1).
# Values that look numeric ("50", "123456") are entered with a leading
# apostrophe so Excel stores them as text (shared strings) instead of
# silently converting them to numbers, matching how the source data was
# authored. The style is then reset to Normal so the cells keep the
# default (unformatted) look, only the underlying value stays text.
$ws.Range("A1").Value = "Armachen Anbessa"
$ws.Range("B1").Value = "Male"
$ws.Range("C1").Value = "'50"
$ws.Range("D1").Value = "'123456"
$ws.Range("E1").Value = "arma@gmail.com"
$ws.Range("F1").Value = "456 Addis"

$ws.Range("C1:D1").Style = "Normal"
